# Burndown Chart Sprint 3 - update "Hrs Trabajadas" (hours worked) for the
# two days that were missing data (rows 13 & 14 => D13, D14).
#
# The "Real Total" column (C) is driven by a shared formula
# (C{n} = SUM(C{n-1} - D{n-1})), so simply writing the new hours into D13/D14
# lets Excel recompute C14:C20 for us - matching the target workbook exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D13").Value = 6
$ws.Range("D14").Value = 7
